$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 183, pushing the existing rows 183:191 down to 184:192.
$ws.Rows.Item(183).Insert()

# The new row 183 repeats the same market/category/quality data as the row
# that used to be 183 (now shifted to 184), except for a new weekly date
# (one week later) and a new Volumen value.
$ws.Range("A183").Value2 = $ws.Range("A184").Value2
$ws.Range("B183").Value2 = $ws.Range("B184").Value2
$ws.Range("C183").Value2 = $ws.Range("C184").Value2
$ws.Range("D183").Value2 = 44516
$ws.Range("E183").Value2 = $ws.Range("E184").Value2
$ws.Range("F183").Value2 = $ws.Range("F184").Value2
$ws.Range("G183").Value2 = $ws.Range("G184").Value2
$ws.Range("H183").Value2 = $ws.Range("H184").Value2
$ws.Range("I183").Value2 = $ws.Range("I184").Value2
$ws.Range("J183").Value2 = 400
$ws.Range("K183").Value2 = $ws.Range("K184").Value2
$ws.Range("L183").Value2 = $ws.Range("L184").Value2
$ws.Range("M183").Value2 = $ws.Range("M184").Value2
$ws.Range("N183").Value2 = $ws.Range("N184").Value2
$ws.Range("O183").Value2 = $ws.Range("O184").Value2
$ws.Range("P183").Value2 = $ws.Range("P184").Value2
$ws.Range("Q183").Value2 = $ws.Range("Q184").Value2
$ws.Range("R183").Value2 = $ws.Range("R184").Value2
